$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.060.57'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.744.70'
$ws.Range("E3").Value = '  +1.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.71%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.42'
$ws.Range("E5").Value = '  -1.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5001'
$ws.Range("E7").Value = '  +8.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3591'
$ws.Range("E8").Value = '  +4.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.62'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07278'
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.065'
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9997'
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.38'
$ws.Range("E13").Value = '  +2.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.980'
$ws.Range("E14").Value = '  +2.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.738.30'
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.882'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.05'
$ws.Range("E17").Value = '  -2.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001039'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06374'
$ws.Range("E19").Value = '  +1.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9999'
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.65'
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.729'
$ws.Range("E22").Value = '  +1.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.107.61'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.33'
$ws.Range("E24").Value = '  +5.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.041'
$ws.Range("E25").Value = '  -4.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.39'
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.99'
$ws.Range("E27").Value = '  +2.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.940.55'
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.177'
$ws.Range("E29").Value = '  +1.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.84'
$ws.Range("E30").Value = '  +1.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.053'
$ws.Range("E31").Value = '  +3.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09582'
$ws.Range("E32").Value = '  +5.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.575'
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.418'
$ws.Range("E34").Value = '  +1.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02205'
$ws.Range("E35").Value = '  +0.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05891'
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.11'
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.435'
$ws.Range("E38").Value = '  +2.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2006'
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.791'
$ws.Range("E40").Value = '  +0.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6048'
$ws.Range("E41").Value = '  +1.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.115'
$ws.Range("E42").Value = '  -1.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.589'
$ws.Range("E43").Value = '  +1.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.88'
$ws.Range("E44").Value = '  +2.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.597'
$ws.Range("E45").Value = '  -0.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5685'
$ws.Range("E46").Value = '  +1.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '120.60'
$ws.Range("E47").Value = '  +1.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.870'
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.107'
$ws.Range("E49").Value = '  +1.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06676'
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.28'
$ws.Range("E51").Value = '  +0.85%  '
